$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bottom border under the year-header row (B3:J3) so the header
# keeps only its top rule, matching the widened table underneath.
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142

# Add the new 2023 column (K) to the right of the existing table, copying
# the formatting of the corresponding cell in the row directly above/left
# so number formats / fonts / fills come along with it.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2023

$ws.Range("B4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2621

$ws.Range("B5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 883

$ws.Range("B6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 1738

# Close off the right edge of the table with a thin right border, since K
# is now the last column.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

# Match the column width used by the rest of the data columns (B:J) for
# the newly occupied columns (K:O, matching the template's pre-existing
# far-right "spare" columns sized the same way).
$ws.Range("K1:O6").EntireColumn.ColumnWidth = 7.83
